$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 150, shifting existing rows 150:162 down to 151:163
$ws.Rows("150:150").Insert()

# Fill in the new row 150 with the new weekly record
$ws.Cells.Item(150, 1).Value = 11
$ws.Cells.Item(150, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(150, 3).Value = "Bíobío"
$ws.Cells.Item(150, 4).Value = 44946
$ws.Cells.Item(150, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(150, 5).Value = 8
$ws.Cells.Item(150, 6).Value = 100112021
$ws.Cells.Item(150, 7).Value = "Ají"
$ws.Cells.Item(150, 8).Value = "Americana (o)"
$ws.Cells.Item(150, 9).Value = "Primera"
$ws.Cells.Item(150, 10).Value = 270
$ws.Cells.Item(150, 11).Value = 20000
$ws.Cells.Item(150, 12).Value = 21000
$ws.Cells.Item(150, 13).Value = 20556
$ws.Cells.Item(150, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(150, 15).Value = "Región Metropolitana"
$ws.Cells.Item(150, 16).Value = 822
$ws.Cells.Item(150, 17).Value = 25
$ws.Cells.Item(150, 18).Value = "Hortaliza"
